# Apply the "add dedicated temporal blocks for thermal units" edit.
#
# Summary of the change (see xml diff):
#  - model_config sheet gets a new column H ("default_temporal_block"),
#    pushing the old H:L (temporal_block..is_active) one column right to I:M.
#  - Two existing rows (old row 6 "rolling_look_ahead" / row 7 "rolling_realisation"
#    false-row) swap places.
#  - Four brand-new rows are appended describing two new temporal blocks:
#    "rolling_look_ahead_ST" and "rolling_look_ahead_ST_nuclear" (each with a
#    RollingHorizon/true row and a Base/false row), followed by one blank
#    spacer row.
#  - node_slack_penalty!C7 formula is repointed from model_config!H3 to
#    model_config!I3 (because temporal_block moved from column H to I).
#  - model_config becomes the active/selected sheet instead of scenarios.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model_config")

# --- 1. Insert the new column H, shifting old H:L to I:M -------------------
$ws.Columns.Item(8).Insert()

# The inserted column picks up formatting from its left neighbour (column G);
# restore the correct style (matches the rest of H:M).
$ws.Range("H1:H9").Style = "40% - Accent6"

# --- 2. Extend the table down to row 13 with correctly-styled blank rows ---
# Old row 9 (blank) already carries the full A:M style pattern once the
# column has been inserted, so replicate its formatting downward.
$ws.Range("A9:M9").Copy() | Out-Null
$ws.Range("A10:M13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- 3. Fix up the values that do not simply follow the column shift -------

# New "default_temporal_block" values for the rows that reference a specific
# temporal block by name (mirrors the corresponding new column I value).
$ws.Range("H3").Value2 = "operation_1year"
$ws.Range("H5").Value2 = "rolling_realisation"

# Old row 6 ("rolling_look_ahead", true) and row 7 ("rolling_realisation",
# false) swap position.
$ws.Range("A6").Value2 = "Base"
$ws.Range("I6").Value2 = "rolling_realisation"
$ws.Range("J6").Value2 = ""
$ws.Range("K6").Value2 = ""
$ws.Range("L6").Value2 = ""
$ws.Range("M6").Value2 = $false

$ws.Range("A7").Value2 = "RollingHorizon"
$ws.Range("H7").Value2 = "rolling_look_ahead"
$ws.Range("I7").Value2 = "rolling_look_ahead"
$ws.Range("J7").Value2 = "1D"
$ws.Range("K7").Value2 = "30D"
$ws.Range("L7").Value2 = "90D"
$ws.Range("M7").Value2 = $true

# New temporal block: rolling_look_ahead_ST
$ws.Range("A9").Value2 = "RollingHorizon"
$ws.Range("B9").Value2 = "operation"
$ws.Range("J9").Value2 = "2D"
$ws.Range("I9").Value2 = "rolling_look_ahead_ST"
$ws.Range("K9").Value2 = "30D"
$ws.Range("L9").Value2 = "90D"
$ws.Range("M9").Value2 = $true

$ws.Range("A10").Value2 = "Base"
$ws.Range("B10").Value2 = "operation"
$ws.Range("I10").Value2 = "rolling_look_ahead_ST"
$ws.Range("M10").Value2 = $false

# New temporal block: rolling_look_ahead_ST_nuclear
$ws.Range("A11").Value2 = "RollingHorizon"
$ws.Range("B11").Value2 = "operation"
$ws.Range("I11").Value2 = "rolling_look_ahead_ST_nuclear"
$ws.Range("J11").Value2 = "30D"
$ws.Range("K11").Value2 = "30D"
$ws.Range("L11").Value2 = "90D"
$ws.Range("M11").Value2 = $true

$ws.Range("A12").Value2 = "Base"
$ws.Range("B12").Value2 = "operation"
$ws.Range("I12").Value2 = "rolling_look_ahead_ST_nuclear"
$ws.Range("M12").Value2 = $false

# Header for the new column (set last so the new shared strings end up in
# the expected order: 2D, rolling_look_ahead_ST, rolling_look_ahead_ST_nuclear,
# default_temporal_block).
$ws.Range("H1").Value2 = "default_temporal_block"

# --- 4. Repoint the node_slack_penalty formula to the moved column ---------
$ws2 = $wb.Worksheets.Item("node_slack_penalty")
$ws2.Range("C7").Formula = "=model_config!I3"

# --- 5. Recalculate and fix up the active sheet / selection ----------------
$excel.CalculateFull()

$ws.Activate()
$ws.Range("R12").Select() | Out-Null
